$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$csv = @"
0,struggle,-0.03444957733154287,0.01032954454422018,0.2995094756285349,-0.0041233403608202,0.0484110713005065,-0.0001527163112768
100,struggle,-0.07984677950541165,-0.1310317392150563,0.2751492535074553,-0.0736092627048492,-0.0381790772080421,0.0797179117798805
200,struggle,-0.09086012840270966,-0.02373796701431304,0.2084361910820008,-1.18019163608551,-4.37715482711792,0.3266601860523224
300,struggle,0.005503296852111615,-0.1945102314154306,0.2710902802646161,0.4163046777248382,-2.031584978103638,1.399950385093689
400,struggle,-0.01686970392863169,-0.689376473426814,-0.4833486676216042,-0.085062988102436,-1.666898488998413,0.3200933933258056
500,struggle,0.7827135920524584,-2.646403878927229,-4.282532811164854,0.2128865420818328,-0.1716531366109848,0.4915938079357147
600,struggle,-1.04310808579127,0.482522944609324,-1.852447708447774,0.4960225820541382,0.8994990587234497,-0.5232060551643372
700,struggle,-2.88519119222959,-0.1530523101488742,0.5539670785268116,1.297935962677002,6.144540786743164,0.5485569834709167
800,struggle,-2.359475612640381,-0.4427683353424072,1.851428031921387,0.395077109336853,5.482057571411133,-0.6198754906654358
900,struggle,-3.12837818264961,1.37843454877535,1.588190158208212,0.6418666839599609,1.02396285533905,0.6579018831253052
1000,struggle,-9.990091621875759,0.5772410134474441,7.240908702214556,0.1765400618314743,1.450194120407104,-0.3078760802745819
1100,struggle,3.089871168136629,1.149596989154817,-7.327601671218906,-0.0943786799907684,-0.3147483170032501,-0.7883216142654419
1200,struggle,-1.657969951629632,-0.6963444749514238,-3.005745251973476,0.4792237877845764,-5.076748371124268,-1.524872422218323
1300,struggle,-2.83851852019628,-0.4919258654117603,-0.4817840854326922,-0.015118914656341,-3.905566930770874,0.6056728959083557
1400,struggle,-1.064056754112247,0.9830632209777808,-6.221896290779103,-0.9416487812995912,-3.088382005691528,-0.4016439020633697
1500,struggle,-0.642586628595987,-0.08318759004274168,-0.6120627919833304,-0.3442225754261017,-2.283719778060913,0.2553416788578033
1600,struggle,-2.159881273905432,-1.468431472778319,-0.03336954116820623,-0.1937969923019409,-0.0872010141611099,0.3181080818176269
1700,struggle,-3.000638425350189,-0.02016200125217193,-1.370798766613007,0.7287622094154358,4.673271656036377,-0.6183483600616455
1800,struggle,-4.049820343653357,1.876313845316566,1.094406366348257,0.5477933883666992,5.67539644241333,-1.208291411399841
1900,struggle,-7.227878630161286,1.69673324127992,3.92280129591624,0.9772316813468932,1.330769896507263,0.8104654550552368
2000,struggle,-2.149182796478235,-1.535205185413374,0.8277400732040205,0.8136724829673767,-0.1959350258111953,-0.294437050819397
2100,struggle,3.039104461669923,-1.047793030738832,-4.546324898799259,-0.1049161031842231,-0.2081523388624191,-0.4677700698375702
2200,struggle,-1.214001258214315,1.53115118543307,-4.196547468503316,-0.6108652353286743,-0.9943358898162842,0.4877758920192718
2300,struggle,-2.15471959114075,0.2546487897634502,-0.9123589172959333,-0.5545129179954529,-0.7066183686256409,-0.1945605874061584
2400,struggle,0.6330445607503264,0.5185447335243235,0.3759170770645095,-0.0233655963093042,-0.0335975885391235,-0.5940664410591125
2500,struggle,-0.2923502524693779,-0.1014364187916107,1.366777941584587,0.1398881375789642,0.0471893399953842,0.531147301197052
2600,struggle,-1.484014511108395,-0.6635967791080455,1.253775984048842,-0.0551305897533893,0.0639881342649459,0.093156948685646
2700,struggle,0.2468122641245533,0.4217223922411623,0.5185143599907569,0.1805106848478317,0.0415388382971286,0.1635591685771942
2800,struggle,-0.3225175539652509,-0.6644023060798641,-0.07338536779085762,-0.1348485052585601,0.1539380401372909,0.1916589736938476
2900,struggle,-0.4933383464813202,-0.1615586355328552,0.263843480497599,-0.0897971913218498,0.1873829066753387,-0.0282525178045034
"@

$lines = $csv -split "`n"
$rowCount = $lines.Count
$colCount = 8

$data = New-Object "object[,]" $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $parts = $lines[$i].Trim() -split ","
    $data[$i,0] = [double]$parts[0]
    $data[$i,1] = [string]$parts[1]
    for ($j = 2; $j -lt $colCount; $j++) {
        $data[$i,$j] = [double]$parts[$j]
    }
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 8))
$rng.Value2 = $data

